$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 16 (shifts rows 16-24 down to 17-25, carrying
# their existing styles/merges along automatically).
$ws.Rows("16:16").Insert()

# Copy the formatting of the (now shifted) row 17 into the newly-blank
# row 16 so the new row reuses the existing "middle data row" cell styles
# instead of Excel fabricating brand-new style records.
$ws.Range("B17:J17").Copy()
$ws.Range("B16:J16").PasteSpecial(-4122)

# --- Row 16: brand-new worker record ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1128054869"
$ws.Range("D16").Value = "CARLO PAOLO CANTILLO QUINTANA"
$ws.Range("E16").Value = "2208"
$ws.Range("F16").Value = 108000
$ws.Range("G16").Value = 2700000

# --- Row 17 (previously row 16: CC/1143386223/JOINER.../2306) now becomes
#     the LUIS EDUARDO record that used to sit in row 18 ---
$ws.Range("C17").Value = "1047456167"
$ws.Range("D17").Value = "LUIS EDUARDO RICARDO MONTERROSA"
$ws.Range("E17").Value = "2312"
$ws.Range("F17").Value = 55811
$ws.Range("G17").Value = 1395280

# --- Row 18 (previously row 17: CC/1143386223/JOINER.../2307) keeps the
#     same worker/period, values unchanged ---
$ws.Range("C18").Value = "1143386223"
$ws.Range("D18").Value = "JOINER CORTECERO MONTERROZA"
$ws.Range("E18").Value = "2307"
$ws.Range("F18").Value = 43307
$ws.Range("G18").Value = 1160000

# --- Row 19 (previously row 18: LUIS EDUARDO/2312) now becomes the
#     JOINER/2306 record that used to sit in row 16 ---
$ws.Range("C19").Value = "1143386223"
$ws.Range("D19").Value = "JOINER CORTECERO MONTERROZA"
$ws.Range("E19").Value = "2306"
$ws.Range("F19").Value = 4640
$ws.Range("G19").Value = 1160000

# --- Header totals: 4 workers/periods now, total mora grew accordingly ---
$ws.Range("E11").Value = 211758
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 4
